$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Authors (E2) with the refreshed authors list (double-space separated)
$ws.Range("E2").Value = "[Jennifer%Lighter%Jennifer.Lighter@nyumc.org%1,  Michael%Phillips%NULL%1,  Sarah%Hochman%NULL%1,  Stephanie%Sterling%NULL%1,  Diane%Johnson%NULL%1,  Fritz%Francois%NULL%1,  Anna%Stachel%NULL%1]"

# Update "Other found locations" (I2) with the new value
$ws.Range("I2").Value = "_PMC"
